# Update inventory exposure simulation with enhanced analysis and presentation plots
#
# The "Price Impact (%)", "Incremental IL (%)" and "IL/Price Impact (%)"
# columns (F, G, H) were recomputed to measure the step from the CURRENT
# price tier to the NEXT (higher) price tier, instead of from the previous
# tier to the current one. The last interior tier (row 23) has no further
# "next" tier to compare against in the refreshed analysis, so it is reset
# to 0 for all three metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 3; $r -le 22; $r++) {
    $rNext = $r + 1

    $priceCur  = $ws.Cells.Item($r, 1).Value2      # column A - Price (USDC)
    $priceNext = $ws.Cells.Item($rNext, 1).Value2
    $valueCur  = $ws.Cells.Item($r, 4).Value2       # column D - Total Value (USDC)
    $valueNext = $ws.Cells.Item($rNext, 4).Value2

    $priceImpact = ($priceCur - $priceNext) / $priceNext * 100
    $incrementalIL = ($valueCur - $valueNext) / $valueNext * 100
    $ilOverImpact = $incrementalIL / $priceImpact * 100

    $ws.Cells.Item($r, 6).Value = $priceImpact      # column F
    $ws.Cells.Item($r, 7).Value = $incrementalIL    # column G
    $ws.Cells.Item($r, 8).Value = $ilOverImpact      # column H
}

# Final tier in this window has no following price point in the refreshed
# analysis window, so it resets to zero.
$ws.Cells.Item(23, 6).Value = 0
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(23, 8).Value = 0
